$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 15.2.1 indicator label cell: drop the period after "15.2.1" (re-authored string)
$ws.Range("B4").Value = "15.2.1 Процесс в переходе на не истощительное ведение лесного хозяйства"

# Selection moved to the edited cell
$ws.Range("B4").Select()
